$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 604, pushing the existing rows (and all
# rows below) down by one. This matches the diff: row 604 becomes the
# new "2026/01/10" entry, and the former rows 604-645 become 605-646.
$ws.Rows.Item(604).Insert()

# Column A holds date-like text ("2026/01/10") that must stay a plain
# string (matching every other row in the sheet) instead of being
# auto-converted into a date serial number. Force text format before
# assigning, then clear the formatting override so the cell ends up
# with no explicit style -- identical to its neighboring rows.
$ws.Range("A604").NumberFormat = "@"
$ws.Range("A604").Value = "2026/01/10"
$ws.Range("A604").ClearFormats()

$ws.Range("B604").Value = "土"
$ws.Range("C604").Value = 10
$ws.Range("D604").Value = 201
